# Edit script: insert "Unique" product-ID columns and new meta columns,
# reorder the shuffled Images arrays, matching the scraper's de-dup pass.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Flipkart Laptop Details"
# Insert a single new column at H ("Unique"), shifting the old H
# column ("Images" list) and everything after it one column right.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Flipkart Laptop Details")

$ws1.Range("H1").EntireColumn.Insert()

$ws1.Range("H1").Value = "Unique"

$ws1.Range("H2").Value = "752c735b4b0a"
$ws1.Range("H3").Value = "3adc9b4b7ed0"

$sheet1Images2 = @'
['https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/x/h/l/-original-imagzk3kcpg6fzxf.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/o/x/n/-original-imagzk3kuu8kduxk.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/d/u/y/-original-imagzk3kukhhjdhq.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/3/j/a/z2-493-thin-and-light-laptop-acer-original-imagr6yjpmhhmpvm.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/l/u/v/z2-493-thin-and-light-laptop-acer-original-imagr6yjuvhjntcb.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/v/t/b/-original-imagvkpfbpgzzqfe.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/a/i/i/-original-imagzk3kfk6bktha.jpeg?q=100&crop=false']
'@
$sheet1Images3 = @'
['https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/f/w/y/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3ggcejgxd.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/i/n/p/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3jtt5ht8g.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/p/j/o/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3vxnbzuzx.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/g/f/r/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3b5edadhw.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/x/u/3/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3kzhphypa.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/z/m/n/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3v9f9uqff.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/i/2/e/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3zty7uwzt.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/q/j/a/a3sp14-31pt-thin-and-light-laptop-acer-original-imah4bj3nvfzgk9g.jpeg?q=100&crop=false', 'https://rukminim2.flixcart.com/image/1664/1664/xif0q/computer/d/p/j/-original-imagwkaymfh66hw3.jpeg?q=100&crop=false']
'@

$ws1.Range("I2").Value = $sheet1Images2
$ws1.Range("I3").Value = $sheet1Images3

# ---------------------------------------------------------------
# Sheet 3: "Amazon Laptop Details"
# Insert four new columns at E ("Meta Title", "Meta Keywords",
# "Meta Description", "Unique"), shifting the old E column
# ("Images" list) and everything after it four columns right.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Amazon Laptop Details")

$ws3.Range("E1:H1").EntireColumn.Insert()

$ws3.Range("E1").Value = "Meta Title"
$ws3.Range("F1").Value = "Meta Keywords"
$ws3.Range("G1").Value = "Meta Description"
$ws3.Range("H1").Value = "Unique"

$ws3.Range("E2").Value = "NA"
$ws3.Range("F2").Value = "NA"
$ws3.Range("G2").Value = "NA"
$ws3.Range("H2").Value = "B0D7PYTTGH"

$sheet3Images2 = @'
['https://m.media-amazon.com/images/I/71p+Rn+JgbL._SL1664_.jpg', 'https://m.media-amazon.com/images/I/61qlqvTsocL._SL1664_.jpg', 'https://m.media-amazon.com/images/I/71W6EmfwI-L._SL1664_.jpg', 'https://m.media-amazon.com/images/I/71MllJzfDfL._SL1664_.jpg', 'https://m.media-amazon.com/images/I/71Mwxs8TX+L._SL1664_.jpg', 'https://m.media-amazon.com/images/I/71V25H7JYBL._SL1664_.jpg']
'@

$ws3.Range("I2").Value = $sheet3Images2
